$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.9
$ws.Range("G2").Value = 4.4
$ws.Range("H2").Value = 2.06
$ws.Range("I2").Value = 2.22
$ws.Range("P2").Value = 1.8
$ws.Range("S2").Value = 3.3
$ws.Range("U2").Value = 2.02
$ws.Range("V2").Value = 1.82
$ws.Range("W2").Value = 1.29
$ws.Range("X2").Value = 15
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 16
$ws.Range("AB2").Value = 17
$ws.Range("AD2").Value = 13.5
$ws.Range("AE2").Value = 30
$ws.Range("AF2").Value = 36
$ws.Range("AG2").Value = 21
$ws.Range("AO2").Value = 23
$ws.Range("F3").Value = 1.44
$ws.Range("G3").Value = 1.6
$ws.Range("J3").Value = 4.2
$ws.Range("L3").Value = 1.36
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 2.04
$ws.Range("O3").Value = 1.28
$ws.Range("Q3").Value = 1.78
$ws.Range("R3").Value = 1.31
$ws.Range("S3").Value = 2.74
$ws.Range("T3").Value = 1.01
$ws.Range("U3").Value = 1.67
$ws.Range("V3").Value = 1.1
$ws.Range("W3").Value = 2.66
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000
$ws.Range("J4").Value = 3.85
$ws.Range("O4").Value = 1.31
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 1.96
$ws.Range("R4").Value = 1.39
$ws.Range("S4").Value = 3.35
$ws.Range("T4").Value = 1.82
$ws.Range("X4").Value = 16
$ws.Range("Y4").Value = 16.5
$ws.Range("AC4").Value = 8.199999999999999
$ws.Range("AN4").Value = 13
$ws.Range("AO4").Value = 1000
$ws.Range("F5").Value = 1.34
$ws.Range("K5").Value = 5.7
$ws.Range("N5").Value = 3.85
$ws.Range("O5").Value = 1.32
$ws.Range("P5").Value = 1.98
$ws.Range("Q5").Value = 1.97
$ws.Range("R5").Value = 1.37
$ws.Range("S5").Value = 3.55
$ws.Range("T5").Value = 2.48
$ws.Range("U5").Value = 1.63
$ws.Range("H6").Value = 2.46
$ws.Range("I6").Value = 2.82
$ws.Range("F7").Value = 1.56
$ws.Range("G7").Value = 1.69
$ws.Range("H7").Value = 5.8
$ws.Range("J7").Value = 4.1
$ws.Range("P7").Value = 2.12
$ws.Range("Q7").Value = 1.73
$ws.Range("S7").Value = 2.86
$ws.Range("T7").Value = 1.81
$ws.Range("U7").Value = 1.99
$ws.Range("AC7").Value = 12.5
$ws.Range("AL7").Value = 1000
$ws.Range("AN7").Value = 9.800000000000001
$ws.Range("S8").Value = 1.74
$ws.Range("AA8").Value = 13
$ws.Range("AD8").Value = 12
$ws.Range("AO8").Value = 2.96
$ws.Range("G10").Value = 1.83
$ws.Range("H10").Value = 4.6
$ws.Range("I10").Value = 4.8
$ws.Range("R10").Value = 1.58
$ws.Range("S10").Value = 2.6
$ws.Range("T10").Value = 1.66
$ws.Range("U10").Value = 2.42
$ws.Range("X10").Value = 23
$ws.Range("Y10").Value = 22
$ws.Range("Z10").Value = 40
$ws.Range("AA10").Value = 120
$ws.Range("AB10").Value = 11.5
$ws.Range("AD10").Value = 19.5
$ws.Range("AE10").Value = 55
$ws.Range("AH10").Value = 17.5
$ws.Range("AI10").Value = 55
$ws.Range("AJ10").Value = 20
$ws.Range("AN10").Value = 8.199999999999999
$ws.Range("AO10").Value = 46
$ws.Range("P11").Value = 2.68
$ws.Range("Q11").Value = 1.56
$ws.Range("R11").Value = 1.67
$ws.Range("T11").Value = 2.8
$ws.Range("X11").Value = 30
$ws.Range("AJ11").Value = 7.6
$ws.Range("F12").Value = 1.42
$ws.Range("G12").Value = 1.44
$ws.Range("H12").Value = 9.199999999999999
$ws.Range("I12").Value = 9.6
$ws.Range("N13").Value = 5
$ws.Range("R13").Value = 1.55
$ws.Range("T13").Value = 1.78
$ws.Range("U13").Value = 2.06
$ws.Range("X13").Value = 28
$ws.Range("Y13").Value = 13
$ws.Range("Z13").Value = 11
$ws.Range("AB13").Value = 32
$ws.Range("AC13").Value = 13
$ws.Range("AD13").Value = 12.5
$ws.Range("AO13").Value = 7.2
$ws.Range("G14").Value = 1.43
$ws.Range("H14").Value = 8.6
$ws.Range("Q14").Value = 1.61
$ws.Range("R14").Value = 1.52
$ws.Range("T14").Value = 1.87
$ws.Range("U14").Value = 1.92
$ws.Range("W14").Value = 3.3
$ws.Range("Y14").Value = 40
$ws.Range("AA14").Value = 370
$ws.Range("AD14").Value = 42
$ws.Range("AJ14").Value = 14.5
